# Tape can create Address and Note objects now.
#
# The exported "tape" header row is renamed from the old Title-Case
# labels (State/Address/City/Zip) to the lower-case field names that
# match the new Address/Note object schema (state/address/city/zipcode).
# The BPO column and all of the data rows are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "state"
$ws.Range("B1").Value = "address"
$ws.Range("C1").Value = "city"
$ws.Range("D1").Value = "zipcode"
$ws.Range("E1").Value = "BPO"

# Matches the selection left behind in the saved workbook.
$ws.Range("B2").Select()
